$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation arrived for "Ciboulette" at
# "Feria Lagunitas de Puerto Montt". Insert a fresh row right above the
# current first data row of that block (row 206), pushing the existing
# rows 206-220 down to 207-221, then populate the new row with the
# latest week's figures.
$ws.Rows(206).Insert()

$ws.Cells.Item(206, 1).Value = 4
$ws.Cells.Item(206, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(206, 3).Value = "Los Lagos"
$ws.Cells.Item(206, 4).Value = 44714
$ws.Cells.Item(206, 5).Value = 10
$ws.Cells.Item(206, 6).Value = 100112039
$ws.Cells.Item(206, 7).Value = "Ciboulette"
$ws.Cells.Item(206, 8).Value = "Sin especificar"
$ws.Cells.Item(206, 9).Value = "Primera"
$ws.Cells.Item(206, 10).Value = 80
$ws.Cells.Item(206, 11).Value = 2500
$ws.Cells.Item(206, 12).Value = 3000
$ws.Cells.Item(206, 13).Value = 2750
$ws.Cells.Item(206, 14).Value = "`$/docena de atados"
$ws.Cells.Item(206, 15).Value = "Región Metropolitana"
$ws.Cells.Item(206, 16).Value = 917
$ws.Cells.Item(206, 17).Value = 3
$ws.Cells.Item(206, 18).Value = "Hortaliza"

$ws.Range("D206").NumberFormat = $ws.Range("D207").NumberFormat
